# Update OT for the "ALBERDI, JUAN BAUTISTA AV. 1091" case (row 5)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 (Caso 6569, CAMPANA 382) was removed entirely; the rows below it
# shift up automatically when the row is deleted.
$ws.Rows.Item(7).Delete()

# After the above deletion, the row that used to be "PATAGONES 2728"
# (originally row 10) is now row 9; remove it too, which shifts the last
# remaining row ("BROWN, ALTE. AV. 881") up into its place.
$ws.Rows.Item(9).Delete()

# Update the OT value for the "ALBERDI, JUAN BAUTISTA AV. 1091" case.
$ws.Range("E5").Value = "ICD30334394"
